$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Задачи")

# New task rows appended below row 38 (tasks 9.1 - 9.3 for invite feature)
$ws.Range("A39").Formula = "=A38+1"
$ws.Range("B39").Value = "9.1 - На элемент Фото игрока добавим иконку пригласить, и мини кнопку ""играем?"". (elementPlayer#2)"
$ws.Range("C39").Value = 42016.468055555553
$ws.Range("D39").Value = 42016.468055555553

$ws.Range("A40").Formula = "=A39+1"
$ws.Range("B40").Value = "9.2 - Функционал отправки приглашения. (LogicInvites.send)"
$ws.Range("C40").Value = 42016.468055555553

$ws.Range("A41").Formula = "=A40+1"
$ws.Range("B41").Value = "9.3 - Функционал принятия приглашения. (LogicInvites.accept)"
$ws.Range("C41").Value = 42016.468055555553

$ws.Range("A42").Formula = "=A41+1"
$ws.Range("A43").Formula = "=A42+1"

# Leave the cursor where the author ended up after typing the new rows
$ws.Range("D39").Select()

Write-Output "edit applied"
